$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: date, hours, task
$ws.Range("B14").Value = 44609
$ws.Range("C14").Value = 4.5
$ws.Range("D14").Value = "Model C. e test dei modelli su AD"

# Row 15: date, hours, task
$ws.Range("B15").Value = 44611
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = "Studio memoria e Inferenza"

# Copy the date cell style (style index 4, numFmtId 14, centered) from B13 onto the new date cells
$ws.Range("B13").Copy()
$ws.Range("B14:B15").PasteSpecial(-4122)  # xlPasteFormats

# Widen column D to fit new text (closest achievable to 55.26953125 char-width units)
$ws.Columns("D").ColumnWidth = 54.5

# Update selection to F5 (matches final cursor position in diff)
$ws.Range("F5").Select()
